$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(19, 2).Value2 = 5145339
$ws.Cells.Item(19, 6).Value2 = 'Hapoel Hadera'
$ws.Cells.Item(19, 7).Value2 = 'MS Ashdod'
$ws.Cells.Item(19, 9).Value2 = 2
$ws.Cells.Item(19, 10).Value2 = 'A'
$ws.Cells.Item(19, 11).Value2 = 2.9
$ws.Cells.Item(19, 13).Value2 = 2.2
$ws.Cells.Item(19, 14).Value2 = 3.6
$ws.Cells.Item(19, 15).Value2 = 3.4
$ws.Cells.Item(19, 16).Value2 = 1.85
$ws.Cells.Item(19, 17).Value2 = 0.5
$ws.Cells.Item(19, 18).Value2 = 1.95
$ws.Cells.Item(19, 19).Value2 = 1.9
$ws.Cells.Item(19, 21).Value2 = 1.975
$ws.Cells.Item(19, 22).Value2 = 1.875
$ws.Cells.Item(19, 24).Value2 = -1
$ws.Cells.Item(19, 25).Value2 = 0.8500000000000001
$ws.Cells.Item(19, 26).Value2 = -1
$ws.Cells.Item(19, 27).Value2 = 0.8999999999999999
$ws.Cells.Item(19, 29).Value2 = 0.875
$ws.Cells.Item(20, 2).Value2 = 5145373
$ws.Cells.Item(20, 6).Value2 = 'Sektzia Nes Tziona'
$ws.Cells.Item(20, 7).Value2 = 'Hapoel Bnei Sakhnin'
$ws.Cells.Item(20, 9).Value2 = 0
$ws.Cells.Item(20, 10).Value2 = 'D'
$ws.Cells.Item(20, 11).Value2 = 3
$ws.Cells.Item(20, 13).Value2 = 2.15
$ws.Cells.Item(20, 14).Value2 = 3.1
$ws.Cells.Item(20, 15).Value2 = 3.25
$ws.Cells.Item(20, 16).Value2 = 2.1
$ws.Cells.Item(20, 17).Value2 = 0.25
$ws.Cells.Item(20, 18).Value2 = 1.975
$ws.Cells.Item(20, 19).Value2 = 1.875
$ws.Cells.Item(20, 21).Value2 = 2
$ws.Cells.Item(20, 22).Value2 = 1.85
$ws.Cells.Item(20, 24).Value2 = 2.25
$ws.Cells.Item(20, 25).Value2 = -1
$ws.Cells.Item(20, 26).Value2 = 0.4875
$ws.Cells.Item(20, 27).Value2 = -0.5
$ws.Cells.Item(20, 29).Value2 = 0.8500000000000001
$ws.Cells.Item(48, 2).Value2 = 5266291
$ws.Cells.Item(48, 6).Value2 = 'Hapoel Jerusalem FC'
$ws.Cells.Item(48, 7).Value2 = 'Sektzia Nes Tziona'
$ws.Cells.Item(48, 9).Value2 = 0
$ws.Cells.Item(48, 10).Value2 = 'H'
$ws.Cells.Item(48, 11).Value2 = 1.6
$ws.Cells.Item(48, 12).Value2 = 3.9
$ws.Cells.Item(48, 13).Value2 = 5
$ws.Cells.Item(48, 14).Value2 = 1.666
$ws.Cells.Item(48, 15).Value2 = 3.75
$ws.Cells.Item(48, 16).Value2 = 4.2
$ws.Cells.Item(48, 17).Value2 = -0.75
$ws.Cells.Item(48, 18).Value2 = 1.975
$ws.Cells.Item(48, 19).Value2 = 1.875
$ws.Cells.Item(48, 21).Value2 = 1.925
$ws.Cells.Item(48, 22).Value2 = 1.925
$ws.Cells.Item(48, 23).Value2 = 0.6659999999999999
$ws.Cells.Item(48, 25).Value2 = -1
$ws.Cells.Item(48, 26).Value2 = 0.4875
$ws.Cells.Item(48, 27).Value2 = -0.5
$ws.Cells.Item(48, 28).Value2 = -1
$ws.Cells.Item(48, 29).Value2 = 0.925
$ws.Cells.Item(49, 2).Value2 = 5145379
$ws.Cells.Item(49, 6).Value2 = 'Hapoel Hadera'
$ws.Cells.Item(49, 7).Value2 = 'Hapoel Bnei Sakhnin'
$ws.Cells.Item(49, 9).Value2 = 3
$ws.Cells.Item(49, 10).Value2 = 'A'
$ws.Cells.Item(49, 11).Value2 = 2.35
$ws.Cells.Item(49, 12).Value2 = 3.25
$ws.Cells.Item(49, 13).Value2 = 2.9
$ws.Cells.Item(49, 14).Value2 = 2.7
$ws.Cells.Item(49, 15).Value2 = 3.25
$ws.Cells.Item(49, 16).Value2 = 2.45
$ws.Cells.Item(49, 17).Value2 = 0
$ws.Cells.Item(49, 18).Value2 = 2.025
$ws.Cells.Item(49, 19).Value2 = 1.825
$ws.Cells.Item(49, 21).Value2 = 2.025
$ws.Cells.Item(49, 22).Value2 = 1.825
$ws.Cells.Item(49, 23).Value2 = -1
$ws.Cells.Item(49, 25).Value2 = 1.45
$ws.Cells.Item(49, 26).Value2 = -1
$ws.Cells.Item(49, 27).Value2 = 0.825
$ws.Cells.Item(49, 28).Value2 = 1.025
$ws.Cells.Item(49, 29).Value2 = -1
$ws.Cells.Item(62, 2).Value2 = 5145342
$ws.Cells.Item(62, 6).Value2 = 'MS Ashdod'
$ws.Cells.Item(62, 7).Value2 = 'Hapoel Bnei Sakhnin'
$ws.Cells.Item(62, 8).Value2 = 1
$ws.Cells.Item(62, 9).Value2 = 1
$ws.Cells.Item(62, 10).Value2 = 'D'
$ws.Cells.Item(62, 12).Value2 = 3.25
$ws.Cells.Item(62, 14).Value2 = 1.85
$ws.Cells.Item(62, 15).Value2 = 3.3
$ws.Cells.Item(62, 16).Value2 = 4.333
$ws.Cells.Item(62, 17).Value2 = -0.5
$ws.Cells.Item(62, 18).Value2 = 1.825
$ws.Cells.Item(62, 19).Value2 = 2.025
$ws.Cells.Item(62, 21).Value2 = 1.925
$ws.Cells.Item(62, 22).Value2 = 1.925
$ws.Cells.Item(62, 23).Value2 = -1
$ws.Cells.Item(62, 24).Value2 = 2.3
$ws.Cells.Item(62, 26).Value2 = -1
$ws.Cells.Item(62, 27).Value2 = 1.025
$ws.Cells.Item(62, 28).Value2 = -0.5
$ws.Cells.Item(62, 29).Value2 = 0.4625
$ws.Cells.Item(63, 2).Value2 = 5266292
$ws.Cells.Item(63, 6).Value2 = 'Maccabi Bnei Raina'
$ws.Cells.Item(63, 7).Value2 = 'Sektzia Nes Tziona'
$ws.Cells.Item(63, 8).Value2 = 3
$ws.Cells.Item(63, 9).Value2 = 2
$ws.Cells.Item(63, 10).Value2 = 'H'
$ws.Cells.Item(63, 12).Value2 = 3.3
$ws.Cells.Item(63, 14).Value2 = 2.25
$ws.Cells.Item(63, 15).Value2 = 3.25
$ws.Cells.Item(63, 16).Value2 = 3.25
$ws.Cells.Item(63, 17).Value2 = -0.25
$ws.Cells.Item(63, 18).Value2 = 1.975
$ws.Cells.Item(63, 19).Value2 = 1.875
$ws.Cells.Item(63, 21).Value2 = 1.95
$ws.Cells.Item(63, 22).Value2 = 1.9
$ws.Cells.Item(63, 23).Value2 = 1.25
$ws.Cells.Item(63, 24).Value2 = -1
$ws.Cells.Item(63, 26).Value2 = 0.9750000000000001
$ws.Cells.Item(63, 27).Value2 = -1
$ws.Cells.Item(63, 28).Value2 = 0.95
$ws.Cells.Item(63, 29).Value2 = -1
$ws.Cells.Item(70, 2).Value2 = 5145343
$ws.Cells.Item(70, 6).Value2 = 'Maccabi Netanya'
$ws.Cells.Item(70, 7).Value2 = 'MS Ashdod'
$ws.Cells.Item(70, 8).Value2 = 0
$ws.Cells.Item(70, 9).Value2 = 2
$ws.Cells.Item(70, 10).Value2 = 'A'
$ws.Cells.Item(70, 11).Value2 = 2.5
$ws.Cells.Item(70, 12).Value2 = 3.2
$ws.Cells.Item(70, 13).Value2 = 2.5
$ws.Cells.Item(70, 14).Value2 = 2.25
$ws.Cells.Item(70, 15).Value2 = 3.1
$ws.Cells.Item(70, 16).Value2 = 2.9
$ws.Cells.Item(70, 17).Value2 = -0.25
$ws.Cells.Item(70, 18).Value2 = 1.975
$ws.Cells.Item(70, 19).Value2 = 1.875
$ws.Cells.Item(70, 20).Value2 = 2.25
$ws.Cells.Item(70, 21).Value2 = 1.9
$ws.Cells.Item(70, 22).Value2 = 1.95
$ws.Cells.Item(70, 23).Value2 = -1
$ws.Cells.Item(70, 25).Value2 = 1.9
$ws.Cells.Item(70, 26).Value2 = -1
$ws.Cells.Item(70, 27).Value2 = 0.875
$ws.Cells.Item(70, 28).Value2 = -0.5
$ws.Cells.Item(70, 29).Value2 = 0.475
$ws.Cells.Item(71, 2).Value2 = 5145154
$ws.Cells.Item(71, 6).Value2 = 'Hapoel Bnei Sakhnin'
$ws.Cells.Item(71, 7).Value2 = 'Beitar Jerusalem'
$ws.Cells.Item(71, 8).Value2 = 2
$ws.Cells.Item(71, 9).Value2 = 0
$ws.Cells.Item(71, 10).Value2 = 'H'
$ws.Cells.Item(71, 11).Value2 = 2.75
$ws.Cells.Item(71, 12).Value2 = 3.25
$ws.Cells.Item(71, 13).Value2 = 2.25
$ws.Cells.Item(71, 14).Value2 = 2.7
$ws.Cells.Item(71, 15).Value2 = 3.2
$ws.Cells.Item(71, 16).Value2 = 2.3
$ws.Cells.Item(71, 17).Value2 = 0
$ws.Cells.Item(71, 18).Value2 = 2.1
$ws.Cells.Item(71, 19).Value2 = 1.775
$ws.Cells.Item(71, 20).Value2 = 2.5
$ws.Cells.Item(71, 21).Value2 = 1.975
$ws.Cells.Item(71, 22).Value2 = 1.875
$ws.Cells.Item(71, 23).Value2 = 1.7
$ws.Cells.Item(71, 25).Value2 = -1
$ws.Cells.Item(71, 26).Value2 = 1.1
$ws.Cells.Item(71, 27).Value2 = -1
$ws.Cells.Item(71, 28).Value2 = -1
$ws.Cells.Item(71, 29).Value2 = 0.875
$ws.Cells.Item(83, 2).Value2 = 6404155
$ws.Cells.Item(83, 6).Value2 = 'Maccabi Haifa'
$ws.Cells.Item(83, 7).Value2 = 'Hapoel Beer Sheva'
$ws.Cells.Item(83, 8).Value2 = 1
$ws.Cells.Item(83, 9).Value2 = 0
$ws.Cells.Item(83, 11).Value2 = 1.6
$ws.Cells.Item(83, 12).Value2 = 3.75
$ws.Cells.Item(83, 13).Value2 = 5
$ws.Cells.Item(83, 14).Value2 = 1.833
$ws.Cells.Item(83, 15).Value2 = 3.3
$ws.Cells.Item(83, 16).Value2 = 4
$ws.Cells.Item(83, 17).Value2 = -0.5
$ws.Cells.Item(83, 18).Value2 = 1.85
$ws.Cells.Item(83, 19).Value2 = 2
$ws.Cells.Item(83, 20).Value2 = 2.5
$ws.Cells.Item(83, 21).Value2 = 1.825
$ws.Cells.Item(83, 22).Value2 = 2.025
$ws.Cells.Item(83, 23).Value2 = 0.833
$ws.Cells.Item(83, 26).Value2 = 0.8500000000000001
$ws.Cells.Item(83, 28).Value2 = -1
$ws.Cells.Item(83, 29).Value2 = 1.025
$ws.Cells.Item(84, 2).Value2 = 6402973
$ws.Cells.Item(84, 6).Value2 = 'MS Ashdod'
$ws.Cells.Item(84, 7).Value2 = 'Maccabi Netanya'
$ws.Cells.Item(84, 8).Value2 = 3
$ws.Cells.Item(84, 9).Value2 = 2
$ws.Cells.Item(84, 11).Value2 = 2.15
$ws.Cells.Item(84, 12).Value2 = 3.3
$ws.Cells.Item(84, 13).Value2 = 3.1
$ws.Cells.Item(84, 14).Value2 = 2.15
$ws.Cells.Item(84, 15).Value2 = 3.2
$ws.Cells.Item(84, 16).Value2 = 3.2
$ws.Cells.Item(84, 17).Value2 = -0.25
$ws.Cells.Item(84, 18).Value2 = 1.925
$ws.Cells.Item(84, 19).Value2 = 1.925
$ws.Cells.Item(84, 20).Value2 = 2.25
$ws.Cells.Item(84, 21).Value2 = 1.8
$ws.Cells.Item(84, 22).Value2 = 2.05
$ws.Cells.Item(84, 23).Value2 = 1.15
$ws.Cells.Item(84, 26).Value2 = 0.925
$ws.Cells.Item(84, 28).Value2 = 0.8
$ws.Cells.Item(84, 29).Value2 = -1
$ws.Cells.Item(85, 2).Value2 = 6404135
$ws.Cells.Item(85, 6).Value2 = 'Hapoel Kiryat Shmona'
$ws.Cells.Item(85, 7).Value2 = 'Hapoel Bnei Sakhnin'
$ws.Cells.Item(85, 8).Value2 = 2
$ws.Cells.Item(85, 9).Value2 = 2
$ws.Cells.Item(85, 11).Value2 = 2.7
$ws.Cells.Item(85, 13).Value2 = 2.375
$ws.Cells.Item(85, 14).Value2 = 2
$ws.Cells.Item(85, 15).Value2 = 3.4
$ws.Cells.Item(85, 16).Value2 = 3.25
$ws.Cells.Item(85, 17).Value2 = -0.25
$ws.Cells.Item(85, 18).Value2 = 1.8
$ws.Cells.Item(85, 19).Value2 = 2.05
$ws.Cells.Item(85, 20).Value2 = 2.25
$ws.Cells.Item(85, 21).Value2 = 1.8
$ws.Cells.Item(85, 22).Value2 = 2.05
$ws.Cells.Item(85, 24).Value2 = 2.4
$ws.Cells.Item(85, 26).Value2 = -0.5
$ws.Cells.Item(85, 27).Value2 = 0.5249999999999999
$ws.Cells.Item(85, 28).Value2 = 0.8
$ws.Cells.Item(85, 29).Value2 = -1
$ws.Cells.Item(86, 2).Value2 = 6402965
$ws.Cells.Item(86, 6).Value2 = 'Maccabi Bnei Raina'
$ws.Cells.Item(86, 7).Value2 = 'Hapoel Hadera'
$ws.Cells.Item(86, 8).Value2 = 1
$ws.Cells.Item(86, 9).Value2 = 1
$ws.Cells.Item(86, 11).Value2 = 2.625
$ws.Cells.Item(86, 13).Value2 = 2.4
$ws.Cells.Item(86, 14).Value2 = 2.6
$ws.Cells.Item(86, 15).Value2 = 3
$ws.Cells.Item(86, 16).Value2 = 2.7
$ws.Cells.Item(86, 17).Value2 = 0
$ws.Cells.Item(86, 18).Value2 = 1.9
$ws.Cells.Item(86, 19).Value2 = 1.95
$ws.Cells.Item(86, 20).Value2 = 2
$ws.Cells.Item(86, 21).Value2 = 1.875
$ws.Cells.Item(86, 22).Value2 = 1.975
$ws.Cells.Item(86, 24).Value2 = 2
$ws.Cells.Item(86, 26).Value2 = 0
$ws.Cells.Item(86, 27).Value2 = -0
$ws.Cells.Item(86, 28).Value2 = 0
$ws.Cells.Item(86, 29).Value2 = -0
$ws.Cells.Item(103, 2).Value2 = 6404130
$ws.Cells.Item(103, 6).Value2 = 'Hapoel Kiryat Shmona'
$ws.Cells.Item(103, 7).Value2 = 'Sektzia Nes Tziona'
$ws.Cells.Item(103, 8).Value2 = 1
$ws.Cells.Item(103, 9).Value2 = 1
$ws.Cells.Item(103, 11).Value2 = 1.75
$ws.Cells.Item(103, 12).Value2 = 3.5
$ws.Cells.Item(103, 13).Value2 = 4
$ws.Cells.Item(103, 14).Value2 = 1.45
$ws.Cells.Item(103, 15).Value2 = 4.2
$ws.Cells.Item(103, 16).Value2 = 5.25
$ws.Cells.Item(103, 17).Value2 = -1.25
$ws.Cells.Item(103, 20).Value2 = 3
$ws.Cells.Item(103, 21).Value2 = 2.025
$ws.Cells.Item(103, 22).Value2 = 1.825
$ws.Cells.Item(103, 24).Value2 = 3.2
$ws.Cells.Item(103, 28).Value2 = -1
$ws.Cells.Item(103, 29).Value2 = 0.825
$ws.Cells.Item(104, 2).Value2 = 6404128
$ws.Cells.Item(104, 6).Value2 = 'Hapoel Haifa'
$ws.Cells.Item(104, 7).Value2 = 'Hapoel Bnei Sakhnin'
$ws.Cells.Item(104, 8).Value2 = 2
$ws.Cells.Item(104, 9).Value2 = 2
$ws.Cells.Item(104, 11).Value2 = 1.8
$ws.Cells.Item(104, 12).Value2 = 3.4
$ws.Cells.Item(104, 13).Value2 = 3.8
$ws.Cells.Item(104, 14).Value2 = 1.95
$ws.Cells.Item(104, 15).Value2 = 3.2
$ws.Cells.Item(104, 16).Value2 = 3.4
$ws.Cells.Item(104, 17).Value2 = -0.5
$ws.Cells.Item(104, 20).Value2 = 2.25
$ws.Cells.Item(104, 21).Value2 = 1.825
$ws.Cells.Item(104, 22).Value2 = 2.025
$ws.Cells.Item(104, 24).Value2 = 2.2
$ws.Cells.Item(104, 28).Value2 = 0.825
$ws.Cells.Item(104, 29).Value2 = -1
$ws.Cells.Item(116, 2).Value2 = 6404124
$ws.Cells.Item(116, 6).Value2 = 'Hapoel TelAviv'
$ws.Cells.Item(116, 7).Value2 = 'Hapoel Kiryat Shmona'
$ws.Cells.Item(116, 8).Value2 = 2
$ws.Cells.Item(116, 9).Value2 = 2
$ws.Cells.Item(116, 10).Value2 = 'D'
$ws.Cells.Item(116, 11).Value2 = 2.05
$ws.Cells.Item(116, 12).Value2 = 3.3
$ws.Cells.Item(116, 13).Value2 = 3.25
$ws.Cells.Item(116, 14).Value2 = 2.45
$ws.Cells.Item(116, 15).Value2 = 3.1
$ws.Cells.Item(116, 16).Value2 = 2.625
$ws.Cells.Item(116, 17).Value2 = 0
$ws.Cells.Item(116, 18).Value2 = 1.85
$ws.Cells.Item(116, 19).Value2 = 2
$ws.Cells.Item(116, 21).Value2 = 2
$ws.Cells.Item(116, 22).Value2 = 1.85
$ws.Cells.Item(116, 24).Value2 = 2.1
$ws.Cells.Item(116, 25).Value2 = -1
$ws.Cells.Item(116, 26).Value2 = 0
$ws.Cells.Item(116, 27).Value2 = -0
$ws.Cells.Item(116, 28).Value2 = 1
$ws.Cells.Item(116, 29).Value2 = -1
$ws.Cells.Item(118, 2).Value2 = 6404123
$ws.Cells.Item(118, 6).Value2 = 'Beitar Jerusalem'
$ws.Cells.Item(118, 7).Value2 = 'Maccabi Bnei Raina'
$ws.Cells.Item(118, 8).Value2 = 0
$ws.Cells.Item(118, 9).Value2 = 1
$ws.Cells.Item(118, 10).Value2 = 'A'
$ws.Cells.Item(118, 11).Value2 = 1.833
$ws.Cells.Item(118, 12).Value2 = 3.5
$ws.Cells.Item(118, 13).Value2 = 3.75
$ws.Cells.Item(118, 14).Value2 = 2.7
$ws.Cells.Item(118, 15).Value2 = 3.3
$ws.Cells.Item(118, 16).Value2 = 2.375
$ws.Cells.Item(118, 17).Value2 = 0.25
$ws.Cells.Item(118, 18).Value2 = 1.75
$ws.Cells.Item(118, 19).Value2 = 2.125
$ws.Cells.Item(118, 21).Value2 = 1.875
$ws.Cells.Item(118, 22).Value2 = 1.975
$ws.Cells.Item(118, 24).Value2 = -1
$ws.Cells.Item(118, 25).Value2 = 1.375
$ws.Cells.Item(118, 26).Value2 = -1
$ws.Cells.Item(118, 27).Value2 = 1.125
$ws.Cells.Item(118, 28).Value2 = -1
$ws.Cells.Item(118, 29).Value2 = 0.9750000000000001
$ws.Cells.Item(127, 2).Value2 = 6404139
$ws.Cells.Item(127, 6).Value2 = 'Maccabi Tel Aviv'
$ws.Cells.Item(127, 7).Value2 = 'Hapoel Beer Sheva'
$ws.Cells.Item(127, 8).Value2 = 3
$ws.Cells.Item(127, 9).Value2 = 0
$ws.Cells.Item(127, 10).Value2 = 'H'
$ws.Cells.Item(127, 11).Value2 = 2.2
$ws.Cells.Item(127, 12).Value2 = 3.3
$ws.Cells.Item(127, 13).Value2 = 2.8
$ws.Cells.Item(127, 14).Value2 = 2
$ws.Cells.Item(127, 15).Value2 = 3.5
$ws.Cells.Item(127, 16).Value2 = 3.3
$ws.Cells.Item(127, 17).Value2 = -0.5
$ws.Cells.Item(127, 18).Value2 = 1.975
$ws.Cells.Item(127, 19).Value2 = 1.875
$ws.Cells.Item(127, 20).Value2 = 2.5
$ws.Cells.Item(127, 23).Value2 = 1
$ws.Cells.Item(127, 25).Value2 = -1
$ws.Cells.Item(127, 26).Value2 = 0.9750000000000001
$ws.Cells.Item(127, 27).Value2 = -1
$ws.Cells.Item(128, 2).Value2 = 6404140
$ws.Cells.Item(128, 6).Value2 = 'Maccabi Netanya'
$ws.Cells.Item(128, 7).Value2 = 'Maccabi Haifa'
$ws.Cells.Item(128, 8).Value2 = 1
$ws.Cells.Item(128, 9).Value2 = 5
$ws.Cells.Item(128, 10).Value2 = 'A'
$ws.Cells.Item(128, 11).Value2 = 4.5
$ws.Cells.Item(128, 12).Value2 = 4.2
$ws.Cells.Item(128, 13).Value2 = 1.533
$ws.Cells.Item(128, 14).Value2 = 6
$ws.Cells.Item(128, 15).Value2 = 5
$ws.Cells.Item(128, 16).Value2 = 1.363
$ws.Cells.Item(128, 17).Value2 = 1.25
$ws.Cells.Item(128, 18).Value2 = 2
$ws.Cells.Item(128, 19).Value2 = 1.85
$ws.Cells.Item(128, 20).Value2 = 3
$ws.Cells.Item(128, 23).Value2 = -1
$ws.Cells.Item(128, 25).Value2 = 0.363
$ws.Cells.Item(128, 26).Value2 = -1
$ws.Cells.Item(128, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(129, 2).Value2 = 6670415
$ws.Cells.Item(129, 6).Value2 = 'Hapoel Beer Sheva'
$ws.Cells.Item(129, 7).Value2 = 'Maccabi Netanya'
$ws.Cells.Item(129, 8).Value2 = 2
$ws.Cells.Item(129, 11).Value2 = 1.5
$ws.Cells.Item(129, 12).Value2 = 4
$ws.Cells.Item(129, 13).Value2 = 5
$ws.Cells.Item(129, 14).Value2 = 1.45
$ws.Cells.Item(129, 15).Value2 = 4.2
$ws.Cells.Item(129, 16).Value2 = 5.25
$ws.Cells.Item(129, 17).Value2 = -1
$ws.Cells.Item(129, 18).Value2 = 1.8
$ws.Cells.Item(129, 19).Value2 = 2.05
$ws.Cells.Item(129, 20).Value2 = 2.75
$ws.Cells.Item(129, 21).Value2 = 1.875
$ws.Cells.Item(129, 22).Value2 = 1.975
$ws.Cells.Item(129, 23).Value2 = 0.45
$ws.Cells.Item(129, 26).Value2 = 0.8
$ws.Cells.Item(129, 28).Value2 = -1
$ws.Cells.Item(129, 29).Value2 = 0.9750000000000001
$ws.Cells.Item(131, 2).Value2 = 6670416
$ws.Cells.Item(131, 6).Value2 = 'Maccabi Haifa'
$ws.Cells.Item(131, 7).Value2 = 'Hapoel Jerusalem FC'
$ws.Cells.Item(131, 8).Value2 = 5
$ws.Cells.Item(131, 11).Value2 = 1.285
$ws.Cells.Item(131, 12).Value2 = 5
$ws.Cells.Item(131, 13).Value2 = 7.5
$ws.Cells.Item(131, 14).Value2 = 1.222
$ws.Cells.Item(131, 15).Value2 = 5.75
$ws.Cells.Item(131, 16).Value2 = 8.5
$ws.Cells.Item(131, 17).Value2 = -1.75
$ws.Cells.Item(131, 18).Value2 = 1.875
$ws.Cells.Item(131, 19).Value2 = 1.975
$ws.Cells.Item(131, 20).Value2 = 3.25
$ws.Cells.Item(131, 21).Value2 = 1.925
$ws.Cells.Item(131, 22).Value2 = 1.925
$ws.Cells.Item(131, 23).Value2 = 0.222
$ws.Cells.Item(131, 26).Value2 = 0.875
$ws.Cells.Item(131, 28).Value2 = 0.925
$ws.Cells.Item(131, 29).Value2 = -1
$ws.Cells.Item(134, 2).Value2 = 6798412
$ws.Cells.Item(134, 6).Value2 = 'Hapoel Beer Sheva'
$ws.Cells.Item(134, 7).Value2 = 'Hapoel Hadera'
$ws.Cells.Item(134, 8).Value2 = 3
$ws.Cells.Item(134, 9).Value2 = 0
$ws.Cells.Item(134, 10).Value2 = 'H'
$ws.Cells.Item(134, 11).Value2 = 1.3
$ws.Cells.Item(134, 12).Value2 = 4.5
$ws.Cells.Item(134, 13).Value2 = 8
$ws.Cells.Item(134, 14).Value2 = 1.333
$ws.Cells.Item(134, 15).Value2 = 4.333
$ws.Cells.Item(134, 16).Value2 = 7
$ws.Cells.Item(134, 17).Value2 = -1.25
$ws.Cells.Item(134, 18).Value2 = 1.825
$ws.Cells.Item(134, 19).Value2 = 2.025
$ws.Cells.Item(134, 20).Value2 = 2.75
$ws.Cells.Item(134, 21).Value2 = 1.95
$ws.Cells.Item(134, 22).Value2 = 1.9
$ws.Cells.Item(134, 23).Value2 = 0.333
$ws.Cells.Item(134, 24).Value2 = -1
$ws.Cells.Item(134, 26).Value2 = 0.825
$ws.Cells.Item(134, 27).Value2 = -1
$ws.Cells.Item(134, 28).Value2 = 0.475
$ws.Cells.Item(134, 29).Value2 = -0.5
$ws.Cells.Item(135, 2).Value2 = 6799822
$ws.Cells.Item(135, 6).Value2 = 'Maccabi Netanya'
$ws.Cells.Item(135, 7).Value2 = 'Maccabi Bnei Raina'
$ws.Cells.Item(135, 8).Value2 = 1
$ws.Cells.Item(135, 9).Value2 = 1
$ws.Cells.Item(135, 10).Value2 = 'D'
$ws.Cells.Item(135, 11).Value2 = 1.85
$ws.Cells.Item(135, 12).Value2 = 3.3
$ws.Cells.Item(135, 13).Value2 = 3.7
$ws.Cells.Item(135, 14).Value2 = 1.65
$ws.Cells.Item(135, 15).Value2 = 3.5
$ws.Cells.Item(135, 16).Value2 = 4.5
$ws.Cells.Item(135, 17).Value2 = -0.75
$ws.Cells.Item(135, 18).Value2 = 1.95
$ws.Cells.Item(135, 19).Value2 = 1.9
$ws.Cells.Item(135, 20).Value2 = 2.5
$ws.Cells.Item(135, 21).Value2 = 2
$ws.Cells.Item(135, 22).Value2 = 1.85
$ws.Cells.Item(135, 23).Value2 = -1
$ws.Cells.Item(135, 24).Value2 = 2.5
$ws.Cells.Item(135, 26).Value2 = -1
$ws.Cells.Item(135, 27).Value2 = 0.8999999999999999
$ws.Cells.Item(135, 28).Value2 = -1
$ws.Cells.Item(135, 29).Value2 = 0.8500000000000001
$ws.Cells.Item(138, 2).Value2 = 6799828
$ws.Cells.Item(138, 6).Value2 = 'Maccabi Bnei Raina'
$ws.Cells.Item(138, 7).Value2 = 'Hapoel Beer Sheva'
$ws.Cells.Item(138, 11).Value2 = 5.5
$ws.Cells.Item(138, 12).Value2 = 4
$ws.Cells.Item(138, 13).Value2 = 1.5
$ws.Cells.Item(138, 14).Value2 = 5.25
$ws.Cells.Item(138, 15).Value2 = 3.8
$ws.Cells.Item(138, 16).Value2 = 1.533
$ws.Cells.Item(138, 17).Value2 = 1
$ws.Cells.Item(138, 18).Value2 = 1.825
$ws.Cells.Item(138, 19).Value2 = 2.025
$ws.Cells.Item(138, 21).Value2 = 1.975
$ws.Cells.Item(138, 22).Value2 = 1.875
$ws.Cells.Item(138, 24).Value2 = 2.8
$ws.Cells.Item(138, 26).Value2 = 0.825
$ws.Cells.Item(138, 27).Value2 = -1
$ws.Cells.Item(138, 29).Value2 = 0.875
$ws.Cells.Item(139, 2).Value2 = 6799830
$ws.Cells.Item(139, 6).Value2 = 'Hapoel Petah Tikva'
$ws.Cells.Item(139, 7).Value2 = 'Hapoel Bnei Sakhnin'
$ws.Cells.Item(139, 11).Value2 = 2.375
$ws.Cells.Item(139, 12).Value2 = 3.2
$ws.Cells.Item(139, 13).Value2 = 2.625
$ws.Cells.Item(139, 14).Value2 = 2.4
$ws.Cells.Item(139, 15).Value2 = 3.1
$ws.Cells.Item(139, 16).Value2 = 2.7
$ws.Cells.Item(139, 17).Value2 = 0
$ws.Cells.Item(139, 18).Value2 = 1.8
$ws.Cells.Item(139, 19).Value2 = 2.05
$ws.Cells.Item(139, 21).Value2 = 2.025
$ws.Cells.Item(139, 22).Value2 = 1.825
$ws.Cells.Item(139, 24).Value2 = 2.1
$ws.Cells.Item(139, 26).Value2 = 0
$ws.Cells.Item(139, 27).Value2 = -0
$ws.Cells.Item(139, 29).Value2 = 0.825
$ws.Cells.Item(140, 2).Value2 = 6799825
$ws.Cells.Item(140, 6).Value2 = 'Hapoel Haifa'
$ws.Cells.Item(140, 7).Value2 = 'Maccabi Petach Tikva'
$ws.Cells.Item(140, 9).Value2 = 2
$ws.Cells.Item(140, 10).Value2 = 'D'
$ws.Cells.Item(140, 11).Value2 = 1.8
$ws.Cells.Item(140, 12).Value2 = 3.25
$ws.Cells.Item(140, 13).Value2 = 4
$ws.Cells.Item(140, 14).Value2 = 1.95
$ws.Cells.Item(140, 15).Value2 = 3.2
$ws.Cells.Item(140, 16).Value2 = 3.4
$ws.Cells.Item(140, 17).Value2 = -0.5
$ws.Cells.Item(140, 18).Value2 = 2.1
$ws.Cells.Item(140, 19).Value2 = 1.775
$ws.Cells.Item(140, 21).Value2 = 2.025
$ws.Cells.Item(140, 22).Value2 = 1.825
$ws.Cells.Item(140, 23).Value2 = -1
$ws.Cells.Item(140, 24).Value2 = 2.2
$ws.Cells.Item(140, 26).Value2 = -1
$ws.Cells.Item(140, 27).Value2 = 0.7749999999999999
$ws.Cells.Item(140, 28).Value2 = 1.025
$ws.Cells.Item(140, 29).Value2 = -1
$ws.Cells.Item(141, 2).Value2 = 6799829
$ws.Cells.Item(141, 6).Value2 = 'Hapoel TelAviv'
$ws.Cells.Item(141, 7).Value2 = 'Maccabi Netanya'
$ws.Cells.Item(141, 9).Value2 = 0
$ws.Cells.Item(141, 10).Value2 = 'H'
$ws.Cells.Item(141, 11).Value2 = 2.4
$ws.Cells.Item(141, 12).Value2 = 3.4
$ws.Cells.Item(141, 13).Value2 = 2.6
$ws.Cells.Item(141, 14).Value2 = 2.625
$ws.Cells.Item(141, 15).Value2 = 3.4
$ws.Cells.Item(141, 16).Value2 = 2.375
$ws.Cells.Item(141, 17).Value2 = 0
$ws.Cells.Item(141, 18).Value2 = 2
$ws.Cells.Item(141, 19).Value2 = 1.85
$ws.Cells.Item(141, 21).Value2 = 2
$ws.Cells.Item(141, 22).Value2 = 1.85
$ws.Cells.Item(141, 23).Value2 = 1.625
$ws.Cells.Item(141, 24).Value2 = -1
$ws.Cells.Item(141, 26).Value2 = 1
$ws.Cells.Item(141, 27).Value2 = -1
$ws.Cells.Item(141, 28).Value2 = -1
$ws.Cells.Item(141, 29).Value2 = 0.8500000000000001
$ws.Cells.Item(146, 2).Value2 = 6799838
$ws.Cells.Item(146, 6).Value2 = 'Hapoel Haifa'
$ws.Cells.Item(146, 7).Value2 = 'MS Ashdod'
$ws.Cells.Item(146, 8).Value2 = 2
$ws.Cells.Item(146, 9).Value2 = 0
$ws.Cells.Item(146, 10).Value2 = 'H'
$ws.Cells.Item(146, 11).Value2 = 2.15
$ws.Cells.Item(146, 12).Value2 = 3
$ws.Cells.Item(146, 13).Value2 = 3.2
$ws.Cells.Item(146, 14).Value2 = 2.15
$ws.Cells.Item(146, 15).Value2 = 3.1
$ws.Cells.Item(146, 16).Value2 = 3.1
$ws.Cells.Item(146, 17).Value2 = -0.25
$ws.Cells.Item(146, 18).Value2 = 2
$ws.Cells.Item(146, 19).Value2 = 1.85
$ws.Cells.Item(146, 20).Value2 = 2.5
$ws.Cells.Item(146, 21).Value2 = 2
$ws.Cells.Item(146, 22).Value2 = 1.85
$ws.Cells.Item(146, 23).Value2 = 1.15
$ws.Cells.Item(146, 24).Value2 = -1
$ws.Cells.Item(146, 26).Value2 = 1
$ws.Cells.Item(146, 27).Value2 = -1
$ws.Cells.Item(146, 29).Value2 = 0.8500000000000001
$ws.Cells.Item(147, 2).Value2 = 6799836
$ws.Cells.Item(147, 6).Value2 = 'Maccabi Tel Aviv'
$ws.Cells.Item(147, 7).Value2 = 'Maccabi Bnei Raina'
$ws.Cells.Item(147, 8).Value2 = 1
$ws.Cells.Item(147, 9).Value2 = 1
$ws.Cells.Item(147, 10).Value2 = 'D'
$ws.Cells.Item(147, 11).Value2 = 1.181
$ws.Cells.Item(147, 12).Value2 = 6
$ws.Cells.Item(147, 13).Value2 = 11
$ws.Cells.Item(147, 14).Value2 = 1.2
$ws.Cells.Item(147, 15).Value2 = 6
$ws.Cells.Item(147, 16).Value2 = 10
$ws.Cells.Item(147, 17).Value2 = -1.75
$ws.Cells.Item(147, 18).Value2 = 1.825
$ws.Cells.Item(147, 19).Value2 = 2.025
$ws.Cells.Item(147, 20).Value2 = 3
$ws.Cells.Item(147, 21).Value2 = 1.85
$ws.Cells.Item(147, 22).Value2 = 2
$ws.Cells.Item(147, 23).Value2 = -1
$ws.Cells.Item(147, 24).Value2 = 5
$ws.Cells.Item(147, 26).Value2 = -1
$ws.Cells.Item(147, 27).Value2 = 1.025
$ws.Cells.Item(147, 29).Value2 = 1
$ws.Cells.Item(154, 2).Value2 = 6799846
$ws.Cells.Item(154, 6).Value2 = 'Hapoel Jerusalem FC'
$ws.Cells.Item(154, 7).Value2 = 'Maccabi Netanya'
$ws.Cells.Item(154, 11).Value2 = 2.8
$ws.Cells.Item(154, 12).Value2 = 3.3
$ws.Cells.Item(154, 13).Value2 = 2.3
$ws.Cells.Item(154, 14).Value2 = 2.5
$ws.Cells.Item(154, 15).Value2 = 3.2
$ws.Cells.Item(154, 16).Value2 = 2.6
$ws.Cells.Item(154, 18).Value2 = 1.9
$ws.Cells.Item(154, 19).Value2 = 1.95
$ws.Cells.Item(154, 21).Value2 = 2.05
$ws.Cells.Item(154, 22).Value2 = 1.8
$ws.Cells.Item(154, 24).Value2 = 2.2
$ws.Cells.Item(154, 29).Value2 = 0.8
$ws.Cells.Item(155, 2).Value2 = 6799841
$ws.Cells.Item(155, 6).Value2 = 'MS Ashdod'
$ws.Cells.Item(155, 7).Value2 = 'Maccabi Petach Tikva'
$ws.Cells.Item(155, 11).Value2 = 2.25
$ws.Cells.Item(155, 12).Value2 = 3.25
$ws.Cells.Item(155, 13).Value2 = 2.75
$ws.Cells.Item(155, 14).Value2 = 2.3
$ws.Cells.Item(155, 15).Value2 = 3.3
$ws.Cells.Item(155, 16).Value2 = 2.7
$ws.Cells.Item(155, 18).Value2 = 1.75
$ws.Cells.Item(155, 19).Value2 = 2.05
$ws.Cells.Item(155, 21).Value2 = 2
$ws.Cells.Item(155, 22).Value2 = 1.85
$ws.Cells.Item(155, 24).Value2 = 2.3
$ws.Cells.Item(155, 29).Value2 = 0.8500000000000001
$ws.Cells.Item(188, 2).Value2 = 7542748
$ws.Cells.Item(188, 6).Value2 = 'MS Ashdod'
$ws.Cells.Item(188, 7).Value2 = 'Hapoel Jerusalem FC'
$ws.Cells.Item(188, 8).Value2 = 2
$ws.Cells.Item(188, 9).Value2 = 0
$ws.Cells.Item(188, 10).Value2 = 'H'
$ws.Cells.Item(188, 11).Value2 = 2.5
$ws.Cells.Item(188, 13).Value2 = 2.625
$ws.Cells.Item(188, 14).Value2 = 2.4
$ws.Cells.Item(188, 15).Value2 = 2.9
$ws.Cells.Item(188, 16).Value2 = 3
$ws.Cells.Item(188, 17).Value2 = -0.25
$ws.Cells.Item(188, 18).Value2 = 2.125
$ws.Cells.Item(188, 19).Value2 = 1.75
$ws.Cells.Item(188, 20).Value2 = 2
$ws.Cells.Item(188, 21).Value2 = 2.05
$ws.Cells.Item(188, 22).Value2 = 1.8
$ws.Cells.Item(188, 23).Value2 = 1.4
$ws.Cells.Item(188, 25).Value2 = -1
$ws.Cells.Item(188, 26).Value2 = 1.125
$ws.Cells.Item(188, 27).Value2 = -1
$ws.Cells.Item(188, 28).Value2 = 0
$ws.Cells.Item(188, 29).Value2 = -0
$ws.Cells.Item(190, 2).Value2 = 7542499
$ws.Cells.Item(190, 6).Value2 = 'Maccabi Petach Tikva'
$ws.Cells.Item(190, 7).Value2 = 'Hapoel Beer Sheva'
$ws.Cells.Item(190, 8).Value2 = 1
$ws.Cells.Item(190, 9).Value2 = 4
$ws.Cells.Item(190, 10).Value2 = 'A'
$ws.Cells.Item(190, 11).Value2 = 2.65
$ws.Cells.Item(190, 13).Value2 = 2.4
$ws.Cells.Item(190, 14).Value2 = 3.2
$ws.Cells.Item(190, 15).Value2 = 3.3
$ws.Cells.Item(190, 16).Value2 = 2.05
$ws.Cells.Item(190, 17).Value2 = 0.25
$ws.Cells.Item(190, 18).Value2 = 2
$ws.Cells.Item(190, 19).Value2 = 1.85
$ws.Cells.Item(190, 20).Value2 = 2.25
$ws.Cells.Item(190, 21).Value2 = 1.85
$ws.Cells.Item(190, 22).Value2 = 2
$ws.Cells.Item(190, 23).Value2 = -1
$ws.Cells.Item(190, 25).Value2 = 1.05
$ws.Cells.Item(190, 26).Value2 = -1
$ws.Cells.Item(190, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(190, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(190, 29).Value2 = -1
$ws.Cells.Item(202, 2).Value2 = 7542719
$ws.Cells.Item(202, 6).Value2 = 'Hapoel Haifa'
$ws.Cells.Item(202, 7).Value2 = 'Maccabi Netanya'
$ws.Cells.Item(202, 8).Value2 = 2
$ws.Cells.Item(202, 10).Value2 = 'H'
$ws.Cells.Item(202, 11).Value2 = 2.6
$ws.Cells.Item(202, 12).Value2 = 3.1
$ws.Cells.Item(202, 13).Value2 = 2.6
$ws.Cells.Item(202, 14).Value2 = 2.9
$ws.Cells.Item(202, 15).Value2 = 3.2
$ws.Cells.Item(202, 16).Value2 = 2.3
$ws.Cells.Item(202, 17).Value2 = 0.25
$ws.Cells.Item(202, 18).Value2 = 1.8
$ws.Cells.Item(202, 19).Value2 = 2.05
$ws.Cells.Item(202, 20).Value2 = 2.5
$ws.Cells.Item(202, 21).Value2 = 2
$ws.Cells.Item(202, 22).Value2 = 1.85
$ws.Cells.Item(202, 23).Value2 = 1.9
$ws.Cells.Item(202, 24).Value2 = -1
$ws.Cells.Item(202, 26).Value2 = 0.8
$ws.Cells.Item(202, 27).Value2 = -1
$ws.Cells.Item(202, 28).Value2 = 1
$ws.Cells.Item(202, 29).Value2 = -1
$ws.Cells.Item(204, 2).Value2 = 7542639
$ws.Cells.Item(204, 6).Value2 = 'Maccabi Bnei Raina'
$ws.Cells.Item(204, 7).Value2 = 'Hapoel Jerusalem FC'
$ws.Cells.Item(204, 8).Value2 = 1
$ws.Cells.Item(204, 10).Value2 = 'D'
$ws.Cells.Item(204, 11).Value2 = 2.5
$ws.Cells.Item(204, 12).Value2 = 3
$ws.Cells.Item(204, 13).Value2 = 2.75
$ws.Cells.Item(204, 14).Value2 = 2.7
$ws.Cells.Item(204, 15).Value2 = 2.8
$ws.Cells.Item(204, 16).Value2 = 2.75
$ws.Cells.Item(204, 17).Value2 = 0
$ws.Cells.Item(204, 18).Value2 = 1.925
$ws.Cells.Item(204, 19).Value2 = 1.925
$ws.Cells.Item(204, 20).Value2 = 2
$ws.Cells.Item(204, 21).Value2 = 2.1
$ws.Cells.Item(204, 22).Value2 = 1.775
$ws.Cells.Item(204, 23).Value2 = -1
$ws.Cells.Item(204, 24).Value2 = 1.8
$ws.Cells.Item(204, 26).Value2 = 0
$ws.Cells.Item(204, 27).Value2 = -0
$ws.Cells.Item(204, 28).Value2 = 0
$ws.Cells.Item(204, 29).Value2 = -0
$ws.Cells.Item(274, 2).Value2 = 6799962
$ws.Cells.Item(274, 6).Value2 = 'MS Ashdod'
$ws.Cells.Item(274, 7).Value2 = 'Hapoel Petah Tikva'
$ws.Cells.Item(274, 8).Value2 = 2
$ws.Cells.Item(274, 11).Value2 = 2.2
$ws.Cells.Item(274, 12).Value2 = 3.1
$ws.Cells.Item(274, 13).Value2 = 3.2
$ws.Cells.Item(274, 14).Value2 = 2.2
$ws.Cells.Item(274, 15).Value2 = 3.1
$ws.Cells.Item(274, 16).Value2 = 3.2
$ws.Cells.Item(274, 17).Value2 = -0.25
$ws.Cells.Item(274, 18).Value2 = 2
$ws.Cells.Item(274, 19).Value2 = 1.85
$ws.Cells.Item(274, 21).Value2 = 2
$ws.Cells.Item(274, 22).Value2 = 1.85
$ws.Cells.Item(274, 23).Value2 = 1.2
$ws.Cells.Item(274, 26).Value2 = 1
$ws.Cells.Item(274, 28).Value2 = -0.5
$ws.Cells.Item(274, 29).Value2 = 0.425
$ws.Cells.Item(275, 2).Value2 = 6799960
$ws.Cells.Item(275, 6).Value2 = 'Maccabi Petach Tikva'
$ws.Cells.Item(275, 7).Value2 = 'Maccabi Bnei Raina'
$ws.Cells.Item(275, 8).Value2 = 1
$ws.Cells.Item(275, 11).Value2 = 2.625
$ws.Cells.Item(275, 12).Value2 = 3.25
$ws.Cells.Item(275, 13).Value2 = 2.5
$ws.Cells.Item(275, 14).Value2 = 2.8
$ws.Cells.Item(275, 15).Value2 = 3.25
$ws.Cells.Item(275, 16).Value2 = 2.375
$ws.Cells.Item(275, 17).Value2 = 0.25
$ws.Cells.Item(275, 18).Value2 = 1.775
$ws.Cells.Item(275, 19).Value2 = 2.1
$ws.Cells.Item(275, 21).Value2 = 1.875
$ws.Cells.Item(275, 22).Value2 = 1.975
$ws.Cells.Item(275, 23).Value2 = 1.8
$ws.Cells.Item(275, 26).Value2 = 0.7749999999999999
$ws.Cells.Item(275, 28).Value2 = -1
$ws.Cells.Item(275, 29).Value2 = 0.9750000000000001
$ws.Cells.Item(279, 18).Value2 = 1.975
$ws.Cells.Item(279, 19).Value2 = 1.875
$ws.Cells.Item(280, 14).Value2 = 1.333
$ws.Cells.Item(280, 15).Value2 = 5.25
$ws.Cells.Item(280, 16).Value2 = 7
$ws.Cells.Item(280, 18).Value2 = 1.85
$ws.Cells.Item(280, 19).Value2 = 2
$ws.Cells.Item(280, 21).Value2 = 1.875
$ws.Cells.Item(280, 22).Value2 = 1.975
$ws.Cells.Item(281, 14).Value2 = 2.05
$ws.Cells.Item(281, 16).Value2 = 3.4
$ws.Cells.Item(281, 18).Value2 = 1.825
$ws.Cells.Item(281, 19).Value2 = 2.025
$ws.Cells.Item(281, 21).Value2 = 1.875
$ws.Cells.Item(281, 22).Value2 = 1.975
$ws.Cells.Item(282, 14).Value2 = 2.1
$ws.Cells.Item(282, 15).Value2 = 3.1
$ws.Cells.Item(282, 16).Value2 = 3.6
$ws.Cells.Item(282, 18).Value2 = 1.825
$ws.Cells.Item(282, 19).Value2 = 2.025
$ws.Cells.Item(283, 14).Value2 = 2.625
$ws.Cells.Item(283, 15).Value2 = 3.25
$ws.Cells.Item(283, 16).Value2 = 2.5
$ws.Cells.Item(283, 21).Value2 = 1.925
$ws.Cells.Item(283, 22).Value2 = 1.925
$ws.Cells.Item(284, 14).Value2 = 1.4
$ws.Cells.Item(284, 15).Value2 = 4.5
$ws.Cells.Item(284, 16).Value2 = 6.5
$ws.Cells.Item(284, 18).Value2 = 1.975
$ws.Cells.Item(284, 19).Value2 = 1.875
$ws.Cells.Item(284, 21).Value2 = 1.925
$ws.Cells.Item(284, 22).Value2 = 1.925
